$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric identifier columns
$ws.Range("A2").Value = 80139019
$ws.Range("B2").Value = 88921

# Taxon id
$ws.Range("E2").Value = 5741

# Species names
$ws.Range("F2").Value = "Tjockfotad fingersvamp"
$ws.Range("G2").Value = "Ramaria flavescens"
$ws.Range("H2").Value = "(Schaeff.) R. H. Petersen"

# Antal (I2) cleared but kept as an (empty-text) cell - use quote-prefix
# trick so the stored type stays Text instead of going fully blank.
$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = "Normal"

# Enhet / Alder-Stadium / Kon / Metod / Bestamningsmetod cells removed entirely
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("AF2").Value = ""

# Location name and coordinates
$ws.Range("P2").Value = "Tveta friluftsgård, 300 m V om, Srm"
$ws.Range("Q2").Value = 648222.682956806
$ws.Range("R2").Value = 6560420.292955686
$ws.Range("S2").Value = 50

# Dates - force text storage (source cells are plain text, not real dates)
# so Excel's date auto-detection doesn't turn them into serial numbers.
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2019-09-27"
$ws.Range("Y2").Style = "Normal"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2019-09-27"
$ws.Range("AA2").Style = "Normal"

# Biotope description added
$ws.Range("AI2").Value = "barrskog"

# Reporter / observer names
$ws.Range("AW2").Value = "Hans Rydberg"
$ws.Range("AX2").Value = "Hans Rydberg"
